# Emerson.xlsx - "Listas sem duplicação de professores"
# Rearranges which weekday/period cell holds a given teacher's class list,
# removing duplicated entries across the row, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "['MCT-3A-Robótica', -, -, -]"

# Row 3 ---------------------------------------------------------------
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "['MCT-3A-Robótica', -, -, -]"
$ws.Range("F3").Value = "-"

# Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "-"

# Row 6 ---------------------------------------------------------------
$ws.Range("B6").Value = "[-, -, -, 'MCT-3A-Robótica']"
$ws.Range("C6").Value = "-"
$ws.Range("F6").Value = "-"

# Row 7 ---------------------------------------------------------------
$ws.Range("B7").Value = "[-, -, -, 'MCT-3A-Robótica']"
$ws.Range("C7").Value = "-"
$ws.Range("F7").Value = "-"

# Row 8 ---------------------------------------------------------------
$ws.Range("B8").Value = "-"

# Row 18 ---------------------------------------------------------------
$ws.Range("C18").Value = "['ELM-2NA-Eletrônica Básica', -]"
$ws.Range("E18").Value = "['ELM-2NA-Eletrônica Básica', -]"

# Row 19 ---------------------------------------------------------------
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "['ELM-2NA-Eletrônica Básica', -]"

# Row 20 ---------------------------------------------------------------
$ws.Range("C20").Value = "-"
$ws.Range("E20").Value = "-"

# Row 21 ---------------------------------------------------------------
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("E21").Value = "['ELM-2NA-Eletrônica Básica', -]"

$wb.Save()
